$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 308, shifting existing rows 308:324 down to 309:325.
$ws.Rows(308).Insert()

# Populate the newly inserted row 308 with the new price record.
$ws.Cells.Item(308, 1).Value  = 10
$ws.Cells.Item(308, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(308, 3).Value  = "La Araucanía"
$ws.Cells.Item(308, 4).Value  = 44706
$ws.Cells.Item(308, 5).Value  = 9
$ws.Cells.Item(308, 6).Value  = 100112009
$ws.Cells.Item(308, 7).Value  = "Acelga"
$ws.Cells.Item(308, 8).Value  = "Sin especificar"
$ws.Cells.Item(308, 9).Value  = "Primera"
$ws.Cells.Item(308, 10).Value = 50
$ws.Cells.Item(308, 11).Value = 10000
$ws.Cells.Item(308, 12).Value = 10000
$ws.Cells.Item(308, 13).Value = 10000
$ws.Cells.Item(308, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(308, 15).Value = "Región Metropolitana"
$ws.Cells.Item(308, 16).Value = 833
$ws.Cells.Item(308, 17).Value = 12
$ws.Cells.Item(308, 18).Value = "Hortaliza"
